{"js": "// Prep for mid term demo\n//\n// 1) \"...analysis of the numerical aspects of IQA.\" ->\n//    \"...analysis of the numerical aspects of IQA algorithms.\"\n//\n// 2) \"This posed a few key considerations. Which IQAs should be used?\" ->\n//    \"This posed a few key considerations. Which IQA algorithms should be\n//    used?\" (i.e. \" algorithm\" is inserted between \"IQA\" and the trailing\n//    \"s\" of \"IQAs\").\n\nconst body = context.document.body;\n\n// --- Edit 1 -----------------------------------------------------------\n// \"aspects of IQA\" (no trailing period) is a unique match in the document,\n// so inserting \" algorithms\" right after it lands exactly between \"IQA\"\n// and the following full stop.\nconst results1 = body.search(\"aspects of IQA\", { matchCase: true });\nresults1.load(\"items\");\nawait context.sync();\n\nif (results1.items.length !== 1) {\n  throw new Error(\n    `Edit 1: expected exactly 1 match for \"aspects of IQA\", found ${results1.items.length}`\n  );\n}\nresults1.items[0].insertText(\" algorithms\", Word.InsertLocation.after);\nawait context.sync();\n\n// --- Edit 2 -------------------------------------------------------------\n// \"Which IQA\" (without the trailing \"s\") is a unique match, so inserting\n// \" algorithm\" right after it turns \"IQAs\" into \"IQA algorithms\".\nconst results2 = body.search(\"Which IQA\", { matchCase: true });\nresults2.load(\"items\");\nawait context.sync();\n\nif (results2.items.length !== 1) {\n  throw new Error(\n    `Edit 2: expected exactly 1 match for \"Which IQA\", found ${results2.items.length}`\n  );\n}\nresults2.items[0].insertText(\" algorithm\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Prep for mid term demo\n#\n# 1) \"...analysis of the numerical aspects of IQA.\" ->\n#    \"...analysis of the numerical aspects of IQA algorithms.\"\n#\n# 2) \"This posed a few key considerations. Which IQAs should be used?\" ->\n#    \"This posed a few key considerations. Which IQA algorithms should be\n#    used?\" (i.e. \" algorithm\" is inserted between \"IQA\" and the trailing\n#    \"s\" of \"IQAs\").\n\n$d = $word.ActiveDocument\n\n# --- Edit 1 --------------------------------------------------------------\n# \"aspects of IQA\" (no trailing period) is a unique match in the document,\n# so inserting \" algorithms\" right after it lands exactly between \"IQA\"\n# and the following full stop.\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Text = \"aspects of IQA\"\n$rng1.Find.MatchCase = $true\n$rng1.Find.MatchWholeWord = $false\n$rng1.Find.Forward = $true\n$found1 = $rng1.Find.Execute()\nif (-not $found1) {\n    throw \"Edit 1: target phrase 'aspects of IQA' not found\"\n}\n$rng1.InsertAfter(\" algorithms\")\n\n# --- Edit 2 ----------------------------------------------------------------\n# \"Which IQA\" (without the trailing \"s\") is a unique match, so inserting\n# \" algorithm\" right after it turns \"IQAs\" into \"IQA algorithms\".\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Text = \"Which IQA\"\n$rng2.Find.MatchCase = $true\n$rng2.Find.MatchWholeWord = $false\n$rng2.Find.Forward = $true\n$found2 = $rng2.Find.Execute()\nif (-not $found2) {\n    throw \"Edit 2: target phrase 'Which IQA' not found\"\n}\n$rng2.InsertAfter(\" algorithm\")\n"}
